# Sprint Backlog Burndown.xlsx
# Commit: "Update Sprint Backlog to reflect tests for the server side"
#   - Adds three new backlog rows to the 'Sprint 2' sheet describing
#     server-side test tasks (Components / Products / Orders), inserted
#     right above the "Estimate Totals" row.
#   - The "Estimate Totals" row (and the two "Bugs from last sprint" /
#     "Require server implementation" rows below it) shift down by 3 rows.
#   - The burndown chart's data series, which pointed at the old totals
#     row, is repointed at the new totals row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 2")

# Insert three blank rows right before the current "Estimate Totals" row
# (row 33). Excel clones the formatting of the row above the insertion
# point into the new rows, and shifts every row at/after 33 down by 3 -
# including the totals row (33 -> 36) and the two bug-tracker rows
# (35/36 -> 38/39).
$ws.Rows("33:35").Insert()

# The freshly-inserted rows inherited column C's body style (fillId 3)
# and columns E:H inherited the "amount remaining" style (fillId 5) from
# the row above (32). Re-stamp them with the formats actually used for
# this kind of row elsewhere in the sheet: C uses the darker highlight
# fill (same as C3), and D:H (the estimate/week columns) use the plain
# number style (same as D3).
$ws.Range("C3").Copy()
$ws.Range("C33:C35").PasteSpecial(-4122) | Out-Null

$ws.Range("D3").Copy()
$ws.Range("E33:H35").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# New backlog rows: "Server" user story, one row per component that now
# needs server-side tests, each estimated/tracked at a flat 3.
$ws.Cells.Item(33, 1).Value = "Server"
$ws.Cells.Item(33, 3).Value = "Implement Tests for the Components Server Side"

$ws.Cells.Item(34, 1).Value = "Server"
$ws.Cells.Item(34, 3).Value = "Implement Tests for the Products Server Side"

$ws.Cells.Item(35, 1).Value = "Server"
$ws.Cells.Item(35, 3).Value = "Implement Tests for the Orders Server Side"

$ws.Range("D33:H35").Value = 3

# Point the selection at the new totals row like the saved workbook does.
$ws.Range("K33").Select()

# Repoint the burndown chart's single series at the relocated totals row
# (D33:H33 -> D36:H36); the cached values are unchanged since the totals
# formulas still sum the same D3:D26-style ranges.
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$ser = $chart.SeriesCollection().Item(1)
$ser.Formula = "=SERIES(,,'Sprint 2'!`$D`$36:`$H`$36,1)"

$wb.Save()
